# Update the "Pais" sheet with the latest COVID snapshot data and re-sort
# by total cases (column B) descending, matching the source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 17:20"

$countries = @(
  "Estados Unidos",
  "Italia",
  "España",
  "China",
  "Alemania",
  "Iran",
  "Francia",
  "Reino Unido",
  "Suiza",
  "Belgica",
  "Paises Bajos",
  "Turquia",
  "Austria",
  "Corea del Sur",
  "Canada",
  "Portugal",
  "Israel",
  "Brasil",
  "Noruega",
  "Australia",
  "Suecia",
  "Chequia",
  "Irlanda",
  "Dinamarca",
  "Malasia",
  "Chile",
  "Rusia",
  "Rumania",
  "Polonia",
  "Filipinas",
  "Luxemburgo",
  "Ecuador",
  "Japon",
  "Pakistan",
  "Tailandia",
  "Arabia Saudita",
  "Indonesia",
  "Finlandia",
  "Sudafrica",
  "India",
  "Grecia",
  "Islandia",
  "Republica Dominicana",
  "Mexico",
  "Panama",
  "Argentina",
  "Peru",
  "Singapur",
  "Serbia",
  "Croacia",
  "Eslovenia",
  "Colombia",
  "Estonia",
  "Hong Kong",
  "Crucero",
  "Catar",
  "Emiratos Arabes Unidos",
  "Egipto",
  "Nueva Zelanda",
  "Irak",
  "Argelia",
  "Marruecos",
  "Barein",
  "Ucrania",
  "Lituania",
  "Armenia",
  "Hungria",
  "Libano",
  "Bosnia y Herzegovina",
  "Bulgaria",
  "Letonia",
  "Principado de Andorra",
  "Eslovaquia",
  "Tunez",
  "Moldavia",
  "Kazajistan",
  "Costa Rica",
  "Republica de Macedonia",
  "Taiwan",
  "Uruguay",
  "Azerbaiyan",
  "Kuwait",
  "Jordania",
  "Republica de Chipre",
  "Reunion",
  "Burkina Faso",
  "Albania",
  "San Marino",
  "Vietnam",
  "Camerun",
  "Oman",
  "Cuba",
  "Senegal",
  "Afganistan",
  "Malta",
  "Islas Feroe",
  "Costa de Marfil",
  "Uzbekistan",
  "Ghana",
  "Bielorrusia",
  "Mauricio",
  "Sri Lanka",
  "Honduras",
  "Nigeria",
  "Venezuela",
  "Brunei",
  "Martinica",
  "Estado de Palestina",
  "Georgia",
  "Camboya",
  "Kirguistan",
  "Bolivia",
  "Guadalupe",
  "Montenegro",
  "Consejo Danes para los Refugiados",
  "Mayotte",
  "Trinidad yTobago",
  "Ruanda",
  "Gibraltar",
  "Liechtenstein",
  "Paraguay",
  "Isla de Man",
  "Kenia",
  "Banglades",
  "Aruba",
  "Monaco",
  "Madagascar",
  "Guayana Francesa",
  "Macao",
  "Puerto Rico",
  "Polinesia Francesa",
  "Jamaica",
  "Guatemala",
  "Zambia",
  "Barbados",
  "Uganda",
  "El Salvador",
  "Guam",
  "Republica de Yibuti",
  "Togo",
  "Mali",
  "Niger",
  "Bermudas",
  "Etiopia",
  "Guinea",
  "Congo",
  "Tanzania",
  "Maldivas",
  "Islas Virgenes de los Estados Unidos",
  "Nueva Caledonia",
  "Gabon",
  "Eritrea",
  "Haiti",
  "San Martin (Parte Francesa)",
  "Bahamas",
  "Birmania",
  "Dominica",
  "Islas Caimanes",
  "Guinea Ecuatorial",
  "Mongolia",
  "Namibia",
  "Curazao",
  "Seychelles",
  "Siria",
  "Groenlandia",
  "Granada",
  "Laos",
  "Suazilandia",
  "Santa Lucia",
  "Libia",
  "Guinea-Bisau",
  "San Cristobal y Nieves",
  "Mozambique",
  "Surinam",
  "Zimbabue",
  "Guyana",
  "Antigua y Barbuda",
  "Republica del Chad",
  "Sudan",
  "Angola",
  "Santa Sede",
  "San Martin (Parte Holandesa)",
  "Benin",
  "San Bartolome",
  "Cabo Verde",
  "Mauritania",
  "Fiyi",
  "Montserrat",
  "Islas Turcas y Caicos",
  "Nepal",
  "Butan",
  "Nicaragua",
  "Gambia",
  "Liberia",
  "Belice",
  "Republica de Africa Central",
  "Islas Virgenes Britanicas",
  "Botsuana",
  "Somalia",
  "Anguila",
  "Timor Oriental",
  "Sierra Leona",
  "Papua Nueva Guinea",
  "San Vicente y las Granadinas",
)

$bvals = @(
  165482,
  101739,
  94417,
  81518,
  68180,
  44605,
  44550,
  25150,
  16186,
  12775,
  12595,
  10827,
  10038,
  9786,
  7474,
  7443,
  4831,
  4681,
  4599,
  4561,
  4435,
  3002,
  2910,
  2815,
  2766,
  2738,
  2337,
  2245,
  2215,
  2084,
  1988,
  1966,
  1953,
  1914,
  1651,
  1563,
  1528,
  1418,
  1326,
  1251,
  1212,
  1135,
  1109,
  1094,
  1075,
  966,
  950,
  926,
  900,
  867,
  802,
  798,
  745,
  714,
  712,
  693,
  664,
  656,
  647,
  630,
  584,
  574,
  567,
  549,
  533,
  532,
  492,
  463,
  411,
  399,
  398,
  370,
  363,
  362,
  353,
  336,
  330,
  329,
  322,
  320,
  298,
  289,
  268,
  262,
  247,
  246,
  243,
  230,
  207,
  193,
  192,
  186,
  175,
  174,
  169,
  169,
  168,
  167,
  152,
  152,
  143,
  142,
  141,
  135,
  135,
  129,
  119,
  117,
  110,
  109,
  107,
  107,
  106,
  105,
  98,
  94,
  85,
  70,
  69,
  65,
  65,
  60,
  59,
  51,
  50,
  49,
  46,
  43,
  41,
  39,
  36,
  36,
  36,
  35,
  34,
  33,
  32,
  32,
  30,
  30,
  28,
  27,
  27,
  25,
  22,
  19,
  19,
  18,
  17,
  16,
  16,
  15,
  15,
  15,
  14,
  14,
  12,
  12,
  12,
  12,
  11,
  11,
  10,
  10,
  10,
  9,
  9,
  9,
  9,
  8,
  8,
  8,
  8,
  8,
  8,
  8,
  7,
  7,
  7,
  7,
  6,
  6,
  6,
  6,
  6,
  6,
  5,
  5,
  5,
  5,
  4,
  4,
  4,
  3,
  3,
  3,
  3,
  3,
  3,
  2,
  1,
  1,
  1,
  1,
)

$cvals = @(
  1694,
  0,
  6461,
  79,
  1295,
  3110,
  0,
  3009,
  264,
  876,
  845,
  0,
  420,
  125,
  26,
  1035,
  136,
  51,
  154,
  101,
  407,
  1,
  0,
  238,
  140,
  289,
  501,
  136,
  160,
  538,
  0,
  0,
  0,
  197,
  127,
  110,
  114,
  66,
  0,
  0,
  0,
  49,
  208,
  101,
  0,
  146,
  0,
  47,
  115,
  77,
  46,
  0,
  30,
  31,
  0,
  0,
  53,
  0,
  58,
  0,
  0,
  18,
  52,
  1,
  42,
  50,
  45,
  17,
  43,
  40,
  22,
  0,
  27,
  0,
  55,
  34,
  0,
  44,
  16,
  0,
  25,
  23,
  0,
  32,
  23,
  0,
  20,
  0,
  3,
  54,
  13,
  16,
  13,
  4,
  13,
  1,
  0,
  18,
  0,
  0,
  15,
  20,
  2,
  4,
  0,
  2,
  26,
  0,
  7,
  2,
  13,
  10,
  0,
  14,
  17,
  12,
  2,
  0,
  0,
  3,
  1,
  11,
  9,
  2,
  0,
  0,
  3,
  0,
  3,
  0,
  0,
  0,
  0,
  0,
  1,
  0,
  2,
  0,
  12,
  0,
  3,
  0,
  0,
  2,
  0,
  0,
  0,
  1,
  0,
  1,
  9,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  2,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
  1,
  0,
  0,
  2,
  1,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
)

$dvals = @(
  5544,
  14620,
  19259,
  76052,
  15824,
  14656,
  7927,
  135,
  1823,
  1696,
  250,
  162,
  1095,
  5408,
  1114,
  43,
  163,
  127,
  13,
  337,
  16,
  25,
  5,
  1,
  537,
  156,
  121,
  220,
  7,
  49,
  80,
  54,
  424,
  76,
  342,
  165,
  81,
  10,
  31,
  102,
  52,
  198,
  5,
  35,
  9,
  228,
  53,
  240,
  42,
  67,
  10,
  15,
  26,
  128,
  603,
  51,
  61,
  150,
  74,
  152,
  37,
  15,
  295,
  8,
  7,
  30,
  37,
  37,
  17,
  17,
  1,
  10,
  3,
  3,
  18,
  22,
  4,
  12,
  39,
  25,
  26,
  73,
  26,
  23,
  1,
  31,
  52,
  13,
  58,
  5,
  34,
  8,
  40,
  5,
  2,
  74,
  6,
  7,
  31,
  47,
  0,
  16,
  3,
  8,
  39,
  45,
  27,
  18,
  21,
  23,
  3,
  0,
  17,
  0,
  2,
  10,
  1,
  0,
  34,
  0,
  1,
  0,
  1,
  25,
  1,
  1,
  0,
  6,
  10,
  1,
  0,
  2,
  10,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
  10,
  2,
  0,
  0,
  1,
  13,
  0,
  0,
  0,
  0,
  1,
  2,
  1,
  0,
  0,
  0,
  1,
  2,
  2,
  2,
  0,
  0,
  2,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  1,
  0,
  0,
  1,
  1,
  0,
  2,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  1,
)

$evals = @(
  156752,
  75528,
  66889,
  2161,
  51674,
  27051,
  33599,
  23226,
  13968,
  10374,
  11306,
  10497,
  8815,
  4216,
  6268,
  7240,
  4651,
  4387,
  4550,
  4205,
  4239,
  2952,
  2851,
  2724,
  2186,
  2570,
  2199,
  1947,
  2176,
  1947,
  1886,
  1850,
  1473,
  1812,
  1299,
  1388,
  1311,
  1391,
  1292,
  1117,
  1114,
  935,
  1053,
  1031,
  1039,
  712,
  873,
  683,
  835,
  794,
  777,
  769,
  715,
  582,
  99,
  641,
  597,
  465,
  572,
  432,
  512,
  526,
  268,
  528,
  518,
  499,
  439,
  414,
  382,
  374,
  397,
  352,
  360,
  349,
  331,
  312,
  324,
  308,
  278,
  294,
  267,
  216,
  237,
  231,
  246,
  203,
  178,
  192,
  149,
  182,
  158,
  172,
  135,
  165,
  167,
  95,
  161,
  158,
  116,
  104,
  139,
  124,
  131,
  125,
  93,
  83,
  90,
  98,
  89,
  86,
  104,
  101,
  85,
  103,
  88,
  83,
  81,
  70,
  35,
  65,
  61,
  60,
  57,
  21,
  49,
  47,
  46,
  37,
  31,
  36,
  36,
  33,
  25,
  35,
  34,
  33,
  32,
  31,
  30,
  28,
  26,
  24,
  17,
  23,
  22,
  19,
  17,
  5,
  17,
  16,
  15,
  15,
  14,
  12,
  13,
  13,
  12,
  11,
  11,
  10,
  9,
  8,
  10,
  8,
  8,
  9,
  9,
  9,
  8,
  8,
  8,
  8,
  8,
  8,
  7,
  7,
  7,
  7,
  4,
  4,
  6,
  6,
  5,
  5,
  5,
  3,
  5,
  5,
  5,
  4,
  4,
  3,
  3,
  3,
  3,
  3,
  3,
  3,
  2,
  2,
  1,
  1,
  1,
  0,
)

$fvals = @(
  3535,
  3981,
  5607,
  528,
  1979,
  3703,
  5056,
  163,
  301,
  1021,
  1053,
  568,
  198,
  55,
  120,
  188,
  83,
  296,
  97,
  28,
  358,
  64,
  103,
  145,
  94,
  14,
  8,
  62,
  50,
  1,
  31,
  58,
  56,
  12,
  23,
  31,
  0,
  56,
  7,
  0,
  72,
  11,
  0,
  1,
  43,
  0,
  49,
  22,
  62,
  32,
  24,
  29,
  13,
  5,
  15,
  6,
  2,
  0,
  2,
  0,
  0,
  1,
  2,
  0,
  27,
  30,
  6,
  7,
  1,
  14,
  3,
  10,
  1,
  10,
  44,
  6,
  7,
  1,
  0,
  9,
  11,
  13,
  5,
  3,
  4,
  0,
  8,
  16,
  3,
  0,
  3,
  3,
  0,
  0,
  2,
  3,
  0,
  8,
  1,
  2,
  1,
  5,
  4,
  0,
  6,
  3,
  15,
  0,
  6,
  1,
  0,
  3,
  10,
  1,
  0,
  3,
  0,
  0,
  0,
  0,
  3,
  0,
  2,
  1,
  0,
  0,
  6,
  0,
  0,
  0,
  2,
  0,
  1,
  0,
  0,
  0,
  5,
  0,
  0,
  0,
  0,
  0,
  0,
  2,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
)

$gvals = @(
  45,
  0,
  553,
  5,
  37,
  141,
  0,
  381,
  36,
  192,
  175,
  0,
  20,
  4,
  3,
  20,
  1,
  4,
  4,
  0,
  34,
  2,
  0,
  13,
  6,
  4,
  8,
  13,
  1,
  10,
  0,
  0,
  0,
  5,
  1,
  2,
  14,
  4,
  0,
  0,
  0,
  0,
  9,
  8,
  0,
  3,
  0,
  0,
  7,
  0,
  4,
  0,
  1,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  0,
  1,
  1,
  2,
  0,
  0,
  0,
  0,
  1,
  2,
  1,
  0,
  2,
  0,
  0,
  1,
  0,
  0,
  1,
  0,
  0,
  2,
  0,
  0,
  0,
  0,
  2,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  1,
  0,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  2,
  0,
  1,
  0,
  1,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
)

$hvals = @(
  3186,
  11591,
  8269,
  3305,
  682,
  2898,
  3024,
  1789,
  395,
  705,
  1039,
  168,
  128,
  162,
  92,
  160,
  17,
  167,
  36,
  19,
  180,
  25,
  54,
  90,
  43,
  12,
  17,
  78,
  32,
  88,
  22,
  62,
  56,
  26,
  10,
  10,
  136,
  17,
  3,
  32,
  46,
  2,
  51,
  28,
  27,
  26,
  24,
  3,
  23,
  6,
  15,
  14,
  4,
  4,
  10,
  1,
  6,
  41,
  1,
  46,
  35,
  33,
  4,
  13,
  8,
  3,
  16,
  12,
  12,
  8,
  0,
  8,
  0,
  10,
  4,
  2,
  2,
  9,
  5,
  1,
  5,
  0,
  5,
  8,
  0,
  12,
  13,
  25,
  0,
  6,
  0,
  6,
  0,
  4,
  0,
  0,
  1,
  2,
  5,
  1,
  4,
  2,
  7,
  2,
  3,
  1,
  2,
  1,
  0,
  0,
  0,
  6,
  4,
  2,
  8,
  1,
  3,
  0,
  0,
  0,
  3,
  0,
  1,
  5,
  0,
  1,
  0,
  0,
  0,
  2,
  0,
  1,
  1,
  0,
  0,
  0,
  0,
  1,
  0,
  1,
  2,
  3,
  0,
  0,
  0,
  0,
  1,
  0,
  0,
  0,
  1,
  0,
  0,
  1,
  0,
  1,
  0,
  1,
  0,
  0,
  0,
  1,
  0,
  2,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  1,
  1,
  0,
  0,
  2,
  2,
  0,
  0,
  0,
  0,
  1,
  1,
  0,
  0,
  0,
  0,
  0,
  1,
  1,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
  0,
)

$startRow = 4
for ($i = 0; $i -lt $countries.Length; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $countries[$i]
  $ws.Cells.Item($r, 2).Value = $bvals[$i]
  $ws.Cells.Item($r, 3).Value = $cvals[$i]
  $ws.Cells.Item($r, 4).Value = $dvals[$i]
  $ws.Cells.Item($r, 5).Value = $evals[$i]
  $ws.Cells.Item($r, 6).Value = $fvals[$i]
  $ws.Cells.Item($r, 7).Value = $gvals[$i]
  $ws.Cells.Item($r, 8).Value = $hvals[$i]
}

Write-Output "Updated $($countries.Length) country rows"
